$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "1.<<include>> Fazer Login" step entirely (old row 7).
# Deleting the whole row shifts rows 8-11 up to 7-10, carrying their
# values AND formatting (so the old "last row" style from row 11 ends
# up on the new last row, row 10).
$ws.Rows("7").Delete()

# Renumber the remaining steps (old 2/3/4/5 -> new 1/2/3/4) now that the
# login step is gone.
$ws.Range("D7").Value = "1. Mostra Menu de opções"
$ws.Range("C8").Value = "2. Seleciona historico de compras"
$ws.Range("D9").Value = "3. Obtém informação"
$ws.Range("D10").Value = "4. Mostra lista de carros comprados"

# Keep the selection in sync with the new bottom-right cell of the table.
$ws.Range("D10").Select()
